$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserManagement")

# New header for the "userAlreadyExists" validation message column
$ws.Range("H3").Value = "userAlreadyExists"

# Copy the formatting pattern of the existing "positive" row (row 4) down to the
# two new negative-scenario rows (7 and 8) before filling in their values.
$ws.Range("A4:U4").Copy()
$ws.Range("A7:U7").PasteSpecial(-4122)
$ws.Range("A4:U4").Copy()
$ws.Range("A8:U8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 7 - negative scenario: username already exists
$ws.Range("A7").Value = "SI_004"
$ws.Range("B7").Value = "continental.automation+sale@gmail.com"
$ws.Range("C7").Value = "Sale_The_Man_1"
$ws.Range("D7").Value = "ALL MY LICENCES"
$ws.Range("E7").Value = "positive"
$ws.Range("F7").Value = "user"
$ws.Range("G7").Value = "USER MANAGEMENT"
$ws.Range("H7").Value = "Username provided already exists"
$ws.Range("I7").Value = "LOGIN"
$ws.Range("J7").Value = "not"
$ws.Range("K7").Value = "4"
$ws.Range("L7").Value = "1"
$ws.Range("M7").Value = "sasa_techn_1"
$ws.Range("N7").Value = "sasa_techn_1"
$ws.Range("O7").Value = "sasa_techn_1@kik.com"
$ws.Range("P7").Value = "successfully created."
$ws.Range("Q7").Value = "Field is required"
$ws.Range("R7").Value = "Field is required"
$ws.Range("S7").Value = "Your information was successfully saved"
$ws.Range("U7").Value = "value"
$ws.Rows.Item(7).RowHeight = 26.25

# Row 8 - negative scenario: username already exists, with a leading space
$ws.Range("A8").Value = "SI_005"
$ws.Range("B8").Value = "continental.automation+sale@gmail.com"
$ws.Range("C8").Value = "Sale_The_Man_1"
$ws.Range("D8").Value = "ALL MY LICENCES"
$ws.Range("E8").Value = "positive"
$ws.Range("F8").Value = "user"
$ws.Range("G8").Value = "USER MANAGEMENT"
$ws.Range("H8").Value = " Username provided already exists"
$ws.Range("I8").Value = "LOGIN"
$ws.Range("J8").Value = "not"
$ws.Range("K8").Value = "4"
$ws.Range("L8").Value = "1"
$ws.Range("M8").Value = "sasa_techn_1"
$ws.Range("N8").Value = "sasa_techn_1"
$ws.Range("O8").Value = "sasa_techn_1@kik.com"
$ws.Range("P8").Value = "successfully created."
$ws.Range("Q8").Value = "Field is required"
$ws.Range("R8").Value = "Field is required"
$ws.Range("S8").Value = "Your information was successfully saved"
$ws.Range("U8").Value = "value"
$ws.Rows.Item(8).RowHeight = 26.25

# Update the saved selection/scroll position for the sheet
$ws.Range("E17").Select()
